$wb = $excel.ActiveWorkbook

# "lethality" sheet gets a new row of test results (row 4) and becomes
# the active/selected sheet with the selection moved to B5.
$ws = $wb.Worksheets.Item("lethality")

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1000
$ws.Range("C4").Value = 1000
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 1000
$ws.Range("F4").Value = 1000

$ws.Activate()
$ws.Range("B5").Select()
